# refactor: invoice & invoice_items scopes / code simplification
#
# The per-provider export sheets recompute their aggregate totals; this
# updates row 2 (the single data row) on each provider sheet with the
# refreshed values.

$wb = $excel.ActiveWorkbook

$data = @{
    "Chassineau" = @(201.39480000000006, 0, 201.39480000000006, 117.3036, 84.0912, 7.050000000000001, 194.34480000000005, 37.488884999999996)
    "Les Vergers de Boulogne" = @(896.1749999999996, 0, 896.1749999999996, 646.1149999999998, 250.06000000000003, 46.32, 849.8549999999998, 170.16225000000014)
    "ABC PEYRAUD" = @(1225.812, 0, 1225.812, 0, 1225.812, 0, 1225.812, 245.1624)
    "SAS Gusto" = @(201.04800000000003, 122.6727, 78.37530000000001, 0, 201.04800000000003, 0, 201.04800000000003, 40.20960000000002)
    "LECA MAREE" = @(155.52, 0, 155.52, 37.14, 118.37999999999998, 0, 155.52, 29.247000000000007)
    "Broko SAS" = @(10102.530000000006, 2621.7299999999996, 7480.800000000001, 0, 10102.530000000006, 12.27, 10090.260000000006, 2020.505999999998)
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $data[$sheetName]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item(2, $i + 1).Value = $values[$i]
    }
}
